# "all microservices deployed, 3 standby"
# Fill in the three remaining personal/passive microservice rows (fixed term
# account, saving account, vip saving account) the same way the other rows
# in the table are already populated: flip the bootstrap checklist columns
# from NO to SI, and add Base URL / PORT / TAG / CALL TO ACTION / GATEWAY
# values + hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param(
        [int]$Row,
        [string]$BaseUrl,
        [string]$Port,
        [string]$Tag,
        [string]$CallToAction,
        [string]$Gateway
    )

    # Columns F..N: bootstrap checklist, switch from "NO" to "SI"
    foreach ($col in @("F", "G", "H", "I", "J", "K", "L", "N")) {
        $ws.Range("$col$Row").Value = "SI"
    }
    # Column M was already "NO" -> stays "NO" (unchanged in the source diff)

    $ws.Range("O$Row").Value = $BaseUrl
    $ws.Range("P$Row").Value = $Port
    $ws.Range("Q$Row").Value = $Tag

    $ws.Hyperlinks.Add($ws.Range("R$Row"), $CallToAction) | Out-Null
    $ws.Range("R$Row").IndentLevel = 2

    $ws.Hyperlinks.Add($ws.Range("S$Row"), $Gateway) | Out-Null
    $ws.Range("S$Row").IndentLevel = 3
}

# Row 18: ms-personal-passive-fixed-term-account
Set-RowData -Row 18 `
    -BaseUrl "https://ms-perpas-fta.azurewebsites.net/" `
    -Port "80:3000" `
    -Tag "ms-perpas-fta" `
    -CallToAction "https://ms-perpas-fta.azurewebsites.net/personal/passive/fixed_term_account/docs/ui" `
    -Gateway "http://gateway-service.eastus.azurecontainer.io/personal/passive/fixed_term_account/docs/ui"

# Row 19: ms-personal-passive-saving-account
Set-RowData -Row 19 `
    -BaseUrl "https://ms-perpas-savingaccount.azurewebsites.net/" `
    -Port "80:3000" `
    -Tag "ms-perpas-savingaccount" `
    -CallToAction "https://ms-perpas-savingaccount.azurewebsites.net/personal/passive/saving_account/docs/ui" `
    -Gateway "http://gateway-service.eastus.azurecontainer.io/personal/passive/saving_account/docs/ui"

# Row 20: ms-personal-passive-vip-saving-account
Set-RowData -Row 20 `
    -BaseUrl "https://ms-perpas-vipsa.azurewebsites.net/" `
    -Port "80:3000" `
    -Tag "ms-perpas-vipsa" `
    -CallToAction "https://ms-perpas-vipsa.azurewebsites.net/personal/passive/vip_saving_account/docs/ui" `
    -Gateway "http://gateway-service.eastus.azurecontainer.io/personal/passive/vip_saving_account/docs/ui"

# Update the view: scroll so column A is leftmost again and select G18
# (best effort - the sandbox may not persist scroll position to the file).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("G18").Select()
